$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.379281
$ws.Range("H2").Value = 22.137843
$ws.Range("I2").Value = 0.2744121884499962
$ws.Range("J2").Value = 0.2744121884499961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 239.0839323333333
$ws.Range("N2").Value = 717.251797
$ws.Range("O2").Value = 0.4086975387666237
$ws.Range("P2").Value = 0.4086975387666237
$ws.Range("Q2").Value = 1764.267519272652
$ws.Range("R2").Value = 15878.40767345387
$ws.Range("S2").Value = 0.1121515860270763
$ws.Range("T2").Value = 0.1121515860270763
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.379281
$ws.Range("H3").Value = 22.137843
$ws.Range("I3").Value = 0.2744121884499962
$ws.Range("J3").Value = 0.2744121884499961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 117.0512696666667
$ws.Range("N3").Value = 351.153809
$ws.Range("O3").Value = 0.2000910950200451
$ws.Range("P3").Value = 0.2000910950200451
$ws.Range("Q3").Value = 863.7542102771097
$ws.Range("R3").Value = 7773.787892493988
$ws.Range("S3").Value = 0.0549074352738067
$ws.Range("T3").Value = 0.05490743527380668
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.379281
$ws.Range("H4").Value = 22.137843
$ws.Range("I4").Value = 0.2744121884499962
$ws.Range("J4").Value = 0.2744121884499961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 171.15883
$ws.Range("N4").Value = 513.47649
$ws.Range("O4").Value = 0.2925842480357353
$ws.Range("P4").Value = 0.2925842480357353
$ws.Range("Q4").Value = 1263.02910220123
$ws.Range("R4").Value = 11367.26191981107
$ws.Range("S4").Value = 0.08028868380948261
$ws.Range("T4").Value = 0.08028868380948259
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.379281
$ws.Range("H5").Value = 22.137843
$ws.Range("I5").Value = 0.2744121884499962
$ws.Range("J5").Value = 0.2744121884499961
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 57.695868
$ws.Range("N5").Value = 173.087604
$ws.Range("O5").Value = 0.09862711817759588
$ws.Range("P5").Value = 0.09862711817759588
$ws.Range("Q5").Value = 425.754022510908
$ws.Range("R5").Value = 3831.786202598172
$ws.Range("S5").Value = 0.02706448333963048
$ws.Range("T5").Value = 0.02706448333963048
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.29805733333333
$ws.Range("H6").Value = 39.894172
$ws.Range("I6").Value = 0.4945128143207339
$ws.Range("J6").Value = 0.4945128143207338
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 239.0839323333333
$ws.Range("N6").Value = 717.251797
$ws.Range("O6").Value = 0.4086975387666237
$ws.Range("P6").Value = 0.4086975387666237
$ws.Range("Q6").Value = 3179.351839647454
$ws.Range("R6").Value = 28614.16655682708
$ws.Range("S6").Value = 0.2021061701014403
$ws.Range("T6").Value = 0.2021061701014403
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.29805733333333
$ws.Range("H7").Value = 39.894172
$ws.Range("I7").Value = 0.4945128143207339
$ws.Range("J7").Value = 0.4945128143207338
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 117.0512696666667
$ws.Range("N7").Value = 351.153809
$ws.Range("O7").Value = 0.2000910950200451
$ws.Range("P7").Value = 0.2000910950200451
$ws.Range("Q7").Value = 1556.554494966794
$ws.Range("R7").Value = 14008.99045470115
$ws.Range("S7").Value = 0.09894761051887987
$ws.Range("T7").Value = 0.09894761051887985
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.29805733333333
$ws.Range("H8").Value = 39.894172
$ws.Range("I8").Value = 0.4945128143207339
$ws.Range("J8").Value = 0.4945128143207338
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 171.15883
$ws.Range("N8").Value = 513.47649
$ws.Range("O8").Value = 0.2925842480357353
$ws.Range("P8").Value = 0.2925842480357353
$ws.Range("Q8").Value = 2276.079934446253
$ws.Range("R8").Value = 20484.71941001628
$ws.Range("S8").Value = 0.1446866599220671
$ws.Range("T8").Value = 0.1446866599220671
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.29805733333333
$ws.Range("H9").Value = 39.894172
$ws.Range("I9").Value = 0.4945128143207339
$ws.Range("J9").Value = 0.4945128143207338
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 57.695868
$ws.Range("N9").Value = 173.087604
$ws.Range("O9").Value = 0.09862711817759588
$ws.Range("P9").Value = 0.09862711817759588
$ws.Range("Q9").Value = 767.2429605604319
$ws.Range("R9").Value = 6905.186645043887
$ws.Range("S9").Value = 0.04877237377834655
$ws.Range("T9").Value = 0.04877237377834654
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.673314666666667
$ws.Range("H10").Value = 5.019944000000001
$ws.Range("I10").Value = 0.06222529534320158
$ws.Range("J10").Value = 0.06222529534320156
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 239.0839323333333
$ws.Range("N10").Value = 717.251797
$ws.Range("O10").Value = 0.4086975387666237
$ws.Range("P10").Value = 0.4086975387666237
$ws.Range("Q10").Value = 400.0626505377076
$ws.Range("R10").Value = 3600.563854839369
$ws.Range("S10").Value = 0.02543132505579274
$ws.Range("T10").Value = 0.02543132505579273
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.673314666666667
$ws.Range("H11").Value = 5.019944000000001
$ws.Range("I11").Value = 0.06222529534320158
$ws.Range("J11").Value = 0.06222529534320156
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 117.0512696666667
$ws.Range("N11").Value = 351.153809
$ws.Range("O11").Value = 0.2000910950200451
$ws.Range("P11").Value = 0.2000910950200451
$ws.Range("Q11").Value = 195.8636062851885
$ws.Range("R11").Value = 1762.772456566696
$ws.Range("S11").Value = 0.01245072748316691
$ws.Range("T11").Value = 0.01245072748316691
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.673314666666667
$ws.Range("H12").Value = 5.019944000000001
$ws.Range("I12").Value = 0.06222529534320158
$ws.Range("J12").Value = 0.06222529534320156
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 171.15883
$ws.Range("N12").Value = 513.47649
$ws.Range("O12").Value = 0.2925842480357353
$ws.Range("P12").Value = 0.2925842480357353
$ws.Range("Q12").Value = 286.4025805685067
$ws.Range("R12").Value = 2577.62322511656
$ws.Range("S12").Value = 0.01820614124679218
$ws.Range("T12").Value = 0.01820614124679217
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.673314666666667
$ws.Range("H13").Value = 5.019944000000001
$ws.Range("I13").Value = 0.06222529534320158
$ws.Range("J13").Value = 0.06222529534320156
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 57.695868
$ws.Range("N13").Value = 173.087604
$ws.Range("O13").Value = 0.09862711817759588
$ws.Range("P13").Value = 0.09862711817759588
$ws.Range("Q13").Value = 96.54334213046401
$ws.Range("R13").Value = 868.8900791741761
$ws.Range("S13").Value = 0.006137101557449749
$ws.Range("T13").Value = 0.006137101557449748
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.540576
$ws.Range("H14").Value = 13.621728
$ws.Range("I14").Value = 0.1688497018860685
$ws.Range("J14").Value = 0.1688497018860685
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 239.0839323333333
$ws.Range("N14").Value = 717.251797
$ws.Range("O14").Value = 0.4086975387666237
$ws.Range("P14").Value = 0.4086975387666237
$ws.Range("Q14").Value = 1085.578765138357
$ws.Range("R14").Value = 9770.208886245215
$ws.Range("S14").Value = 0.06900845758231434
$ws.Range("T14").Value = 0.06900845758231433
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.540576
$ws.Range("H15").Value = 13.621728
$ws.Range("I15").Value = 0.1688497018860685
$ws.Range("J15").Value = 0.1688497018860685
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 117.0512696666667
$ws.Range("N15").Value = 351.153809
$ws.Range("O15").Value = 0.2000910950200451
$ws.Range("P15").Value = 0.2000910950200451
$ws.Range("Q15").Value = 531.4801858179947
$ws.Range("R15").Value = 4783.321672361952
$ws.Range("S15").Value = 0.03378532174419162
$ws.Range("T15").Value = 0.03378532174419161
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.540576
$ws.Range("H16").Value = 13.621728
$ws.Range("I16").Value = 0.1688497018860685
$ws.Range("J16").Value = 0.1688497018860685
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 171.15883
$ws.Range("N16").Value = 513.47649
$ws.Range("O16").Value = 0.2925842480357353
$ws.Range("P16").Value = 0.2925842480357353
$ws.Range("Q16").Value = 777.1596756860799
$ws.Range("R16").Value = 6994.437081174719
$ws.Range("S16").Value = 0.04940276305739343
$ws.Range("T16").Value = 0.04940276305739343
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.540576
$ws.Range("H17").Value = 13.621728
$ws.Range("I17").Value = 0.1688497018860685
$ws.Range("J17").Value = 0.1688497018860685
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 57.695868
$ws.Range("N17").Value = 173.087604
$ws.Range("O17").Value = 0.09862711817759588
$ws.Range("P17").Value = 0.09862711817759588
$ws.Range("Q17").Value = 261.972473539968
$ws.Range("R17").Value = 2357.752261859712
$ws.Range("S17").Value = 0.01665315950216911
$ws.Range("T17").Value = 0.01665315950216911
